$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New test case row (row 4) - values entered in authoring order so new
# shared-string entries land at the same indices as the authored workbook
$ws.Range("J4").Value = "Philip Revak"
$ws.Range("A4").Value = "CreateAccount_1"
$ws.Range("B4").Value = "User"
$ws.Range("C4").Value = "P1"
$ws.Range("E4").Value = "Website is open in a suported browser"
$ws.Range("D4").Value = "Verify that the proccess of account creation functions properly"
$ws.Range("G4").Value = "Account is successfully created"
$ws.Range("F4").Value = "1. Navigate to the create account page`n2.Enter a valid username in the username box.`n3. Enter a valid password in the password box.`n4. Click the create account button."
$ws.Range("H4").Value = "Not yet tested"
$ws.Range("I4").Value = "Not Executed"

# Copy formatting from row 3 to row 4
$ws.Range("A3:J3").Copy()
$ws.Range("A4:J4").PasteSpecial(-4122)  # xlPasteFormats

# Set a distinct gray fill for the Status cell (Not Executed) instead of the
# green Pass fill: xlThemeColorLight1 (white background theme), darkened.
$ws.Range("I4").Interior.ThemeColor = 2
$ws.Range("I4").Interior.TintAndShade = -0.249977111117893

# Row height for the wrapped text
$ws.Rows.Item(4).RowHeight = 69.599999999999994

# Adjust column width for Status column
$ws.Columns.Item(9).ColumnWidth = 11.88671875

# Update selection
$ws.Range("B4").Select()
